$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.289.74"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "1.551.04"
$ws.Range("E3").Value = "  -1.62%  "
$ws.Range("E4").Value = "  -0.19%  "
$c = $ws.Range("D5")
$c.Value = "'209.95"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.68%  "
$c = $ws.Range("D6")
$c.Value = "'0.480"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("E7").Value = "  -0.11%  "
$c = $ws.Range("D8")
$c.Value = "'23.77"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("E9").Value = "  -2.15%  "
$ws.Range("E10").Value = "  -1.63%  "
$c = $ws.Range("D11")
$c.Value = "'0.0889"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").Value = "1.773.38"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").Value = "1.550.86"
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("D14").Value = "28.304.43"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("E16").Value = "  -1.89%  "
$c = $ws.Range("D17")
$c.Value = "'60.80"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.52%  "
$c = $ws.Range("D18")
$c.Value = "'227.47"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("E20").Value = "  -2.69%  "
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("E23").Value = "  -2.77%  "
$ws.Range("E24").Value = "  -1.33%  "
$c = $ws.Range("D25")
$c.Value = "'150.75"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("E29").Value = "  -3.54%  "
$ws.Range("E30").Value = "  -2.92%  "
$ws.Range("E31").Value = "  -4.94%  "
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("D33").Value = "1.387.42"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("E35").Value = "  +2.63%  "
$c = $ws.Range("D36")
$c.Value = "'1.48"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -4.20%  "
$ws.Range("E37").Value = "  -1.10%  "
$c = $ws.Range("D38")
$c.Value = "'2.57"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.93%  "
$ws.Range("E39").Value = "  -3.22%  "
$c = $ws.Range("D40")
$c.Value = "'0.513"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.23%  "
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").Value = "  -2.23%  "
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("E45").Value = "  -2.35%  "
$c = $ws.Range("D46")
$c.Value = "'61.85"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.12%  "
$ws.Range("D47").Value = "1.685.55"
$ws.Range("E47").Value = "  -1.52%  "
$c = $ws.Range("D48")
$c.Value = "'0.897"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -6.74%  "
$c = $ws.Range("D49")
$c.Value = "'85.46"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.30%  "
$c = $ws.Range("D50")
$c.Value = "'42.62"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +7.37%  "
$ws.Range("E51").Value = "  +0.31%  "
